$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final, chronologically-sorted list of saved-ad URLs (rows 2..33).
# Rows 2-28 are rendered as live hyperlinks; rows 29-33 are plain text
# (matches the target workbook exactly).
$urls = @(
    "https://www.autotrader.co.uk/car-details/202302204461549?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202306058163586?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202405139652007?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202407011312408?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202409204270639?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202410084967287?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202410295706788?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202410295712702?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202411015835018?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202411146282450?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202411156319287?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202412036942163?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202412077088745?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202412317639076?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202501158046351?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202501188152626?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202501308577450?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202502018634758?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202502038698884?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202502078873648?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202502199256509?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202502219353755?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202502229382614?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202502259476344?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202503029658728?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202503039693826?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202503019623450?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202503049742784?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202503079858704?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202503059791560?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202503049732509?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ",
    "https://www.autotrader.co.uk/car-details/202503059777856?fromSavedAds=true&advertising-location=at_cars&sort=relevance&postcode=CB58TJ"
)

$hyperlinkCount = 27   # rows 2..28 get real hyperlinks; the rest stay plain text

# Capture the existing "Hyperlink" cell style (s=2 in the original file)
# from A2 before anything else touches it, so every hyperlink row below can
# be re-stamped with that *same* style rather than a freshly-minted one.
$hyperlinkStyle = $ws.Range("A2").Style()

# --- wipe out the old hyperlinks + old data rows (keep row 1, the header) ---
$ws.Hyperlinks.Delete()
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count()
if ($lastRow -lt ($urls.Length + 1)) { $lastRow = $urls.Length + 1 }
$ws.Range("A2:A" + $lastRow).ClearContents()

# Pass 1: write every cell value and create the live hyperlinks first, while
# every target cell still carries the plain "Normal" style. Doing the Add()
# calls before any re-styling means they all mint/reuse the very same
# implicit hyperlink xf instead of branching into several near-duplicates.
for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Range("A" + $row)
    $url = $urls[$i]

    if ($i -lt $hyperlinkCount) {
        $ws.Hyperlinks.Add($cell, $url)
    } else {
        $cell.Value = $url
    }
}

# Pass 2: re-stamp the hyperlink rows with the workbook's original
# "Hyperlink" style (s=2) so no extra cellXfs entries stick around.
for ($i = 0; $i -lt $hyperlinkCount; $i++) {
    $row = $i + 2
    $ws.Range("A" + $row).Style = $hyperlinkStyle
}
